# Apply cryptos list update (price/volume refresh + row shift for new ARBITRUM entry)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking text value need to be forced to Text format
# first, otherwise Excel will silently convert strings like "1.000" or "2.740" into
# plain numbers (1, 2.74) and drop the significant trailing zeros / multi-dot format
# used by this sheet to display prices as plain text.
$numericLookingCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D15",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update cell values to match the refreshed data
$ws.Range("D2").Value = "28.500.02"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.909.23"
$ws.Range("E3").Value = "  +5.30%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "313.76"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "0.3951"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "0.09792"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "1.161"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").Value = "41.57"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "6.551"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "21.14"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("D14").Value = "1.909.96"
$ws.Range("D15").Value = "7.581"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "0.00001143"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "93.87"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "0.06656"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "18.09"
$ws.Range("E20").Value = "  +5.65%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "6.301"
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("D23").Value = "28.552.41"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +3.03%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("D26").Value = "2.740"
$ws.Range("E26").Value = "  +14.48%  "
$ws.Range("D27").Value = "2.130.13"
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("D28").Value = "21.37"
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("D29").Value = "159.39"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "129.02"
$ws.Range("D31").Value = "1.105"
$ws.Range("E31").Value = "  +6.52%  "
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "5.715"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "3.635"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "9.939"
$ws.Range("E35").Value = "  +10.37%  "
$ws.Range("D36").Value = "0.06797"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "0.02446"
$ws.Range("E37").Value = "  +5.06%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.277"
$ws.Range("E38").Value = "  +9.89%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2235"
$ws.Range("E39").Value = "  +4.66%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "5.099"
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("D41").Value = "11.65"
$ws.Range("E41").Value = "  +3.62%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.6447"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.193"
$ws.Range("E43").Value = "  +3.92%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.70"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6118"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "1.281"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.662"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.049"
$ws.Range("E49").Value = "  +5.80%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "125.44"
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.211"
$ws.Range("E51").Value = "  +2.63%  "

# Restore default cell style on the cells we temporarily forced to Text format,
# so their formatting matches the rest of the (unstyled) data cells.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
